# Fruta / hortaliza, semanal
# Insert one new week's worth of price records (3 rows: Especial/Primera/Segunda,
# Región Metropolitana, fecha 2022-01-04) right above the existing row 5, pushing
# the rest of the table down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new rows, shifting rows 5..27 down to 8..30.
$ws.Rows("5:7").Insert()

# New row 5: Especial
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5, 3).Value = "Los Lagos"
$ws.Cells.Item(5, 4).Value = 44565
$ws.Cells.Item(5, 5).Value = 10
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100103
$ws.Cells.Item(5, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(5, 9).Value = 100103003
$ws.Cells.Item(5, 10).Value = "Damasco"
$ws.Cells.Item(5, 11).Value = "Castle Brite"
$ws.Cells.Item(5, 12).Value = "Especial"
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 20000
$ws.Cells.Item(5, 15).Value = 20000
$ws.Cells.Item(5, 16).Value = 20000
$ws.Cells.Item(5, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(5, 18).Value = "Región Metropolitana"
$ws.Cells.Item(5, 19).Value = 1111
$ws.Cells.Item(5, 20).Value = 18

# New row 6: Primera
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(6, 3).Value = "Los Lagos"
$ws.Cells.Item(6, 4).Value = 44565
$ws.Cells.Item(6, 5).Value = 10
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100103
$ws.Cells.Item(6, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(6, 9).Value = 100103003
$ws.Cells.Item(6, 10).Value = "Damasco"
$ws.Cells.Item(6, 11).Value = "Castle Brite"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 200
$ws.Cells.Item(6, 14).Value = 18000
$ws.Cells.Item(6, 15).Value = 18000
$ws.Cells.Item(6, 16).Value = 18000
$ws.Cells.Item(6, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(6, 18).Value = "Región Metropolitana"
$ws.Cells.Item(6, 19).Value = 1000
$ws.Cells.Item(6, 20).Value = 18

# New row 7: Segunda
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(7, 3).Value = "Los Lagos"
$ws.Cells.Item(7, 4).Value = 44565
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100103
$ws.Cells.Item(7, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(7, 9).Value = 100103003
$ws.Cells.Item(7, 10).Value = "Damasco"
$ws.Cells.Item(7, 11).Value = "Castle Brite"
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 16000
$ws.Cells.Item(7, 15).Value = 16000
$ws.Cells.Item(7, 16).Value = 16000
$ws.Cells.Item(7, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(7, 18).Value = "Región Metropolitana"
$ws.Cells.Item(7, 19).Value = 889
$ws.Cells.Item(7, 20).Value = 18
